# Generate Report for Handback
# Inserts a new handed-back file (81a93002-d663-4bd8-aee3-b91f461471fa) into
# the Overview / zh-cn / de-de report sheets, positioned (alphabetically) just
# before the existing a6f858de-5453-4d7b-b184-b74f39f2e80d entry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Overview sheet (columns A:G)
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Shift the existing a6f858de row (row 3) down to row 4, duplicating styles.
$ov.Rows.Item(3).Insert()

$ov.Range("A3").Value = "81a93002-d663-4bd8-aee3-b91f461471fa.md"
$ov.Range("B3").Value = "e2e\81a93002-d663-4bd8-aee3-b91f461471fa.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = "2016-09-07 09:09:33"

# Resize the Overview table to include the new row.
$ovTable = $ov.ListObjects.Item("Overview")
$ovTable.Resize($ov.Range("A1:G4"))

# Rebuild hyperlinks on column B in the new, final row order.
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc393c325ba88e3bfddcaf4e1c0ba720ff5124df/e2e/9ead3d87-ee17-42fd-897d-997a1e0c5bb2.md", "", "", "e2e\9ead3d87-ee17-42fd-897d-997a1e0c5bb2.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af71edb88178191bba057e0571c32f53328668d/e2e/81a93002-d663-4bd8-aee3-b91f461471fa.md", "", "", "e2e\81a93002-d663-4bd8-aee3-b91f461471fa.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f7a8e9988904b8350d270cdb9350d881ec82e50/e2e/a6f858de-5453-4d7b-b184-b74f39f2e80d.md", "", "", "e2e\a6f858de-5453-4d7b-b184-b74f39f2e80d.md") | Out-Null

# ---------------------------------------------------------------------------
# 2) zh-cn sheet (columns A:P)
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Rows.Item(3).Insert()

$zh.Range("A3").Value = "81a93002-d663-4bd8-aee3-b91f461471fa.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "'True"
$zh.Range("G3").Value = "81a93002-d663-4bd8-aee3-b91f461471fa.af71edb88178191bba057e0571c32f53328668de.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-07 09:09:25"
$zh.Range("I3").Value = "81a93002-d663-4bd8-aee3-b91f461471fa.md"
$zh.Range("J3").Value = "81a93002-d663-4bd8-aee3-b91f461471fa.af71edb88178191bba057e0571c32f53328668de.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-07 09:09:48"
$zh.Range("L3").Value = ""
$zh.Range("M3").Value = "'True"
$zh.Range("N3").Value = ""
$zh.Range("O3").Value = "'False"
$zh.Range("P3").Value = ""

$zhTable = $zh.ListObjects.Item("zh-cn")
$zhTable.Resize($zh.Range("A1:P4"))

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc393c325ba88e3bfddcaf4e1c0ba720ff5124df/e2e/9ead3d87-ee17-42fd-897d-997a1e0c5bb2.md", "", "", "9ead3d87-ee17-42fd-897d-997a1e0c5bb2.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bf1a43f10b3f8c1dce23d5d216d5dead6b21a21a/e2e/9ead3d87-ee17-42fd-897d-997a1e0c5bb2.md", "", "", "9ead3d87-ee17-42fd-897d-997a1e0c5bb2.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af71edb88178191bba057e0571c32f53328668d/e2e/81a93002-d663-4bd8-aee3-b91f461471fa.md", "", "", "81a93002-d663-4bd8-aee3-b91f461471fa.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/af71edb88178191bba057e0571c32f53328668d/e2e/81a93002-d663-4bd8-aee3-b91f461471fa.md", "", "", "81a93002-d663-4bd8-aee3-b91f461471fa.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f7a8e9988904b8350d270cdb9350d881ec82e50/e2e/a6f858de-5453-4d7b-b184-b74f39f2e80d.md", "", "", "a6f858de-5453-4d7b-b184-b74f39f2e80d.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/51756bb49a268e6193d98e48489bb8d6f6def4f9/e2e/a6f858de-5453-4d7b-b184-b74f39f2e80d.md", "", "", "a6f858de-5453-4d7b-b184-b74f39f2e80d.md") | Out-Null

# ---------------------------------------------------------------------------
# 3) de-de sheet (columns A:P)
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Rows.Item(3).Insert()

$de.Range("A3").Value = "81a93002-d663-4bd8-aee3-b91f461471fa.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "'True"
$de.Range("G3").Value = "81a93002-d663-4bd8-aee3-b91f461471fa.af71edb88178191bba057e0571c32f53328668de.de-de.xlf"
$de.Range("H3").Value = "2016-09-07 09:09:33"
$de.Range("I3").Value = "81a93002-d663-4bd8-aee3-b91f461471fa.md"
$de.Range("J3").Value = "81a93002-d663-4bd8-aee3-b91f461471fa.af71edb88178191bba057e0571c32f53328668de.de-de.xlf"
$de.Range("K3").Value = "2016-09-07 09:09:56"
$de.Range("L3").Value = ""
$de.Range("M3").Value = "'True"
$de.Range("N3").Value = ""
$de.Range("O3").Value = "'False"
$de.Range("P3").Value = ""

$deTable = $de.ListObjects.Item("de-de")
$deTable.Resize($de.Range("A1:P4"))

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cc393c325ba88e3bfddcaf4e1c0ba720ff5124df/e2e/9ead3d87-ee17-42fd-897d-997a1e0c5bb2.md", "", "", "9ead3d87-ee17-42fd-897d-997a1e0c5bb2.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a3b5474a578337e65cb7b60bba332d85316344d5/e2e/9ead3d87-ee17-42fd-897d-997a1e0c5bb2.md", "", "", "9ead3d87-ee17-42fd-897d-997a1e0c5bb2.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/af71edb88178191bba057e0571c32f53328668d/e2e/81a93002-d663-4bd8-aee3-b91f461471fa.md", "", "", "81a93002-d663-4bd8-aee3-b91f461471fa.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/af71edb88178191bba057e0571c32f53328668d/e2e/81a93002-d663-4bd8-aee3-b91f461471fa.md", "", "", "81a93002-d663-4bd8-aee3-b91f461471fa.md") | Out-Null
$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2f7a8e9988904b8350d270cdb9350d881ec82e50/e2e/a6f858de-5453-4d7b-b184-b74f39f2e80d.md", "", "", "a6f858de-5453-4d7b-b184-b74f39f2e80d.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dcf808733d98fd8a2e81a072de36930a9d375a7c/e2e/a6f858de-5453-4d7b-b184-b74f39f2e80d.md", "", "", "a6f858de-5453-4d7b-b184-b74f39f2e80d.md") | Out-Null

Write-Output "Done"
